$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('B2').Value = 'Jaipur Rajasthan performed spectacularly NEET National Eligibility cum Entrance Tests exam year Out 79,057 appeared NEET 45,336 candidates cleared test 182 Jaipur alone four top 50 state.Prince Chaudhary small village Barmer The 17-year-old became famous overnight securing fifth rank country Son medical store owner Prince wants pursue medical course All India Institute Medical Science Delhi become neuro physician `` happy shy teenager told NDTV And though studied Hindi medium school till class 10 two years intense coaching career institute Kota gave confidence take exam English scoring 686 720.Not Prince one every five students polished exam skills coaching centres Kota claims one teachers institute `` Coaching institutes create awareness able help students focused preparation said Naveen Maheshwari heads one sought coaching institutes Kota The city famous cotton sarees seen mushrooming coaching centres thriving NEET boom There around 1.5 lakh students 150 centres city.googletag.cmd.push function googletag.display `` adslotNativeVideo Students mostly enroll institutes class 11 spend two years studying exactly crack exam says Manan Bhardwaj Jaipur got India rank 150 `` It matter much study matters study Manan told NDTV Comments We seeing child Barmer cleared exam also got India rank 5 unheard till years ago Last year saw Rupa Yadav child bride married 8 clear exam So certainly awareness coaching help children back ground Mr Maheshwari said.In coaching capital 150 institutes annual turnover Rs 1500 crore say reports And ''s tuition fees lakhs aspirants come Uttar Pradesh Madhya Pradesh Chhattisgarh Bihar parts northern India rent rooms duration courses There booming PG accommodation business private houses dormitories facilities'
$ws.Range('C2').Value = 'news'
$ws.Range('B3').Value = 'Netflix made name back really good TV shows continues spend billions dollars year build portfolio original series prominently showcases homepage attract eyeballs self-produced efforts By tends background popular critically-acclaimed series n''t make arguably better anything Netflix made history.There ''s need highlight casual viewers may miss Netflix ''s ''ve taken upon gather best TV shows streaming Netflix To prepare list used aggregate ratings Rotten Tomatoes Metacritic IMDb draw shortlist picked top 50 Here best TV shows Netflix India sorted alphabetically This list updated every month worthy additions movies removed service bookmark page keep checking in.The 100 Best Movies Netflix India A Series Unfortunate Events 2017 – Present Three resourceful orphans – inventor reader baby sharp teeth – must outsmart evil distant relative ''s sizeable fortune trying uncover mystery behind secret society parents involved Based Lemony Snicket ''s popular series children ''s books Alias Grace 2017 Margaret Atwood ''s 1996 novel name 19th-century Canadian woman convicted double murder becomes subject criminal psychologist – profession n''t exist name – adapted screen six-part miniseries American Crime Story 2016 – Present A true crime anthology series prolific producer Ryan Murphy follows well-known events dominated US media trial former sportsman O.J Simpson assassination fashion designer Gianni Versace Better Call Saul 2015 – Present This spin-off prequel Breaking Bad follows small-time lawyer Bob Odenkirk tendencies con artist transforms morally-challenged criminal lawyer knew Saul Goodman Black Mirror 2011 – Present Charlie Brooker ''s anthology series consisting standalone episodes – means ever-changing cast addition new settings storylines – explores unanticipated consequences new technologies often dark satirical ways BoJack Horseman 2014 – Present Set world humans anthropomorphic animals live alongside washed-up sitcom star plans comeback fame help ghostwriter ex-girlfriend ''s also agent freeloading roommate dealing rival ''s dating ghostwriter Breaking Bad 2008 – 2013 Diagnosed lung cancer struggling high school chemistry teacher Bryan Cranston decides get business making selling meth secure family ''s financial future help former student Aaron Paul Broadchurch 2013 – 2017 A look violent crimes affect small seaside town Britain eyes two investigators David Tennant Olivia Colman impact media attention suspicion grief close-knit community Brooklyn Nine-Nine 2013 – Present The lives group detectives fictional New York precinct – ensemble featuring Any Samberg AndreBraugher Terry Crews – get sitcom treatment The Office co-creator Michael Schur Crazy Ex-Girlfriend 2015 – Present A successful young female lawyer Rachel Bloom also co-creator suffers depression anxiety uproots life New York moves suburb California find love happiness The Crown 2016 – Present A look life Britain ''s longest-ruling monarch Queen Elizabeth II wedding 1947 present day including political rivalries romances world-changing events reign Daredevil 2015 – Present A blind man fights injustice attorney day masked vigilante night present-day New York struggling deal Catholic guilt arises actions Dear White People 2017 – Present Set fictional predominantly white Ivy League college Justin Simien turns successful satirical film broader series black students dealing social injustice figuring The End F***ing World 2017 – Present Two teenage outsiders – budding psychopath rebel hungry adventure – embark road trip search latter ''s real father become involved series increasingly violent events The Expanse 2015 – Present Hundreds years future mankind colonised Solar System brink war ''s crew different origins – Earth Mars Asteroid Belt – expose greatest conspiracy Fargo 2014 – Present The Coen brothers acclaimed film name inspiration black comedy/crime anthology series featuring quirky characters across different eras dealing deception intrigue murder amidst cold American Midwest Fauda 2015 – Present A former expert Israeli agent comes retirement hunt Palestinian militant show cleverly following sides conflict drawing real world events Fullmetal Alchemist Brotherhood 2009 – 2010 This direct adaptation famous manga centres two brothers looking Philosopher ''s Stone restore bodies failed experimental attempt resurrect mother goes awry But ''re ones stone GLOW 2017 – Present An exploration real-life 1980s women ''s professional wrestling – Gorgeous Ladies Of Wresting GLOW – focusing personal professional lives fictionalised Hollywood misfits involved Godless 2017 In 19th century murderous outlaw notorious gang hunt ex-protégé betrayed find conflict town gives refuge whose residents mainly women A seven-episode miniseries The Good Place 2016 – Present The Office co-creator Michael Schur spun sitcom web afterlife series following woman Kristen Bell ''s mistakenly assigned Heaven-like utopia tries better person hide The Good Wife 2009 – 2016 After humiliating sex corruption scandal puts husband behind bars wife – former state ''s attorney – must return work provide family battling unwanted spotlight Hannibal 2013 – 2015 Cancelled three seasons cult horror hit explores relationship forensic psychiatrist Mads Mikkelsen patient young FBI criminal profiler uncanny ability empathise serial killers Happy Valley 2014 – Present A British police sergeant dealing personal problems investigates cases small West Yorkshire town Two seasons six episodes possibly way House Cards 2013 – Present A betrayed US Congressman works equally conniving wife climb political ladders Washington anything needed – manipulation deceit even murder – achieve goals The Inbetweeners 2008 – 2010 A coming-of-age sitcom whose success led two movies following four British teenagers final year school series misadventures involving uncaring school staff male bonding failed sexual encounters Jane Virgin 2014 – Present Gina Rodriguez Annihilation stars devout Catholic working young Latina virgin becomes pregnant accidental artificial insemination rom-com satire end fifth season 2019 Jessica Jones 2015 – Present Suffering PTSD super-powered woman rebuilds life private investigator New York forced battle past demons every step way wants look past Lost 2004 – 2010 The survivors plane crash must work together survive mysterious tropical island somewhere South Pacific Ocean battling supernatural sci-fi elements Faltered midway recovered ended disappointing fashion Luther 2010 – Present Idris Elba stars dedicated brilliant British detective tries keep grip personal life dealing psychological fallouts crimes ''s tasked solve Mad Men 2007 – 2015 Set 1960s New York slow-burn drama offers peek inside fictional ad agency focusing one extremely talented executives Jon Hamm ''s bored simple personal life Master None 2015 – Present Loosely based Aziz Ansari ''s life experiences 30-year-old struggling actor tries get life order personal professional fronts affected Indian ancestry Mindhunter 2017 – Present In late 1970s two FBI agents push superiors expand research criminal science involves getting close personal imprisoned serial killers understand think Monty Python ''s Flying Circus 1969 – 1974 The influential British comedy group got start surreal sketch series targeted life island intellectual fashion humour wide-ranging unique gave birth term `` Pythonesque Narcos 2015 – Present A gripping look violent powerful drug cartels Colombia including infamous Pablo Escobar thecorroborative efforts various law enforcement whose job bring Orange Is New Black 2013 – Present The lives incarcerated women minimum-security federal prison upstate New York including normally law-abiding privileged woman ''s sentenced decade-old crime Orphan Black 2013 – 2017 A con artist Tatiana Maslany assumes identity woman committed suicide looked like pulled conspiracy learns ''s clone Outlander 2014 – Present Diana Gabaldon ''s best-selling books married nurse World War II ''s transported back time 1743 finds caught Jacobite risings two different men Peaky Blinders 2013 – Present The exploits Shelby crime family Birmingham England two World Wars elements borrowed 19th-century gang name legend goes used sew razor blades caps Peep Show 2003 – 2015 The lives two different dysfunctional twenty-somethings sharing flat South London day-to-day farcical antics made cult favourite Rick Morty 2013 – Present An adult animated sci-fi series cynical alcoholic mad scientist Rick weak-willed easily influenced grandson Morty go series misadventures across dimensions universes Sherlock 2010 – Present Benedict Cumberbatch Martin Freeman play famous detective doctor sidekick modern-day adaptation Sir Arthur Conan Doyle ''s stories Terrific compelling first years though derailed fourth season Star Trek 1966 – 1969 Gene Rodenberry ''s original series adventures starship 2260s led Kirk Spock McCoy n''t age well visually storytelling endures set aside major dip quality third season Star Trek The Next Generation 1987 – 1994 Rodenberry took vision even far-off future – 24th century – focusing new generation Starfleet officers new ship recognisable name Stranger Things 2016 – Present A nostalgia-heavy love letter 80s set suburban small town secret scientific experiments paranormal supernatural times human subjects starts wreak havoc normal lives residents That ''70s Show 1997 – 2006 The comedy fueled Hollywood careers stars Mila Kunis Ashton Kutcher group six high school friends going usual teenage problems trying figure identities Trollhunters 2016 – Present Guillermo del Toro turns love monster stories animated series ages following teenage boy stumbles onto mystic amulet must protect world humans trolls Unbreakable Kimmy Schmidt 2015 – Present Rescued doomsday cult 15 years young woman armed positive attitude decides start new life New York help gay wannabe-Broadway actor roommate street-wise landlady depressed out-of-touch socialite The Walking Dead 2010 – Present Based popular comic series horror drama set post-apocalyptic future survivors search safe world overrun zombies Hit peak fifth season never recovered Wentworth 2013 – Present Locked awaiting trial attempted murder husband woman adjusts life Australian prison rises ranks Gripping till season 4 post lead actress left div class=''compare-product-widget /div'
$ws.Range('C3').Value = 'gadgets'
$ws.Range('B4').Value = 'New Delhi TIMES Higher Education THE World University Ranking released rankings educational institutes Asia two Indian institutes made way top 50 Asian institutes The two institutes made way top 50 Indian Institute Science 29th place Indian Institute Technology Bombay 44th rank Other institutes made way top 100 Indian Institute Technology Kharagpur Indian Institute Technology Roorkee Indian Institute Technology Kanpur Indian Institute Technology Delhi The ranks released THE World University Ranking Indian Institute Science IISc world rank 29Indian Institute Technology IIT Bombay world rank 44Indian Institute Technology IIT Kharagpur world rank 60Indian Institute Technology IIT Roorkee world rank 65Indian Institute Technology IIT Kanpur world rank 81Indian Institute Technology IIT Delhi world rank 86Among made top 200 Tezpur University 100th place Indian Institute Technology Madras 103rd place Indian Institute Technology Guwahati 112th place Panjab University 114th place National Institute Technology Rourkela 126th place Jadavpur University 127th place Indian School Mines 141st place University Delhi 144th place Aligarh Muslim University 158th place Savitribai Phule Pune University 188th place IIT-BHU 194th place CommentsA total 17 institutes India made way top 200 Asian Universities Though representation India increased list several universities fallen rankings comparison last year Both IISc IIT Bombay dropped two places 29th 44th position respectively IIT Madras witnessed biggest decline year fallen last year ''s 41st place 103rd pace year googletag.cmd.push function googletag.display `` adslotNativeVideo Click Education News'
$ws.Range('B5').Value = 'New Delhi A differently abled Indian teacher little-known institution Uttarakhand beat tough competition around world named today among top 50 contenders USD 1 million global award Pradeep Negi economics social science computer science teacher Government Inter College BHEL Ranipur Haridwar shortlisted annual Varkey Foundation Global Teacher Prize 2018 30,000 nominations applications 173 countries around world `` It unbelievable proud chosen one top 50 teachers It motivated teachers community said Mr Negi `` handicapped n''t stopped succeeding glad inspire disabled people India developed lot innovative work ICT poor children glad recognised said Mr Negi physically disabled age two due polio.Read Maggie McDonnell Teacher From Canada Wins Global Teacher Prize 2017Negi chosen fighting spirit suffered bullying school disability overcoming adversities use new technology tools make economics social subjects compelling students.googletag.cmd.push function googletag.display `` adslotNativeVideo Congratulations Pradeep Negi reaching final 50 hope story inspires looking enter teaching profession shines spotlight incredible work teachers world every day said Sunny Varkey founder Varkey Foundation Global Teacher Prize `` We intend keep momentum going journey continues return teachers rightful position one respected professions society said.Mr Negi helped train 1,200 teachers state.The Global Teacher Prize set recognise one exceptional teacher every year made outstanding contribution profession.The top 50 shortlisted teachers narrowed 10 finalists committee results announced February 2018.The winner chosen final 10 finalists Global Teacher Prize Academy.All 10 finalists invited Dubai Award ceremony Global Education Skills Forum GESF March 18 2018 winner announced live stage CommentsClick Education News Except headline story edited NDTV staff published syndicated feed'
$ws.Range('B6').Value = 'New Delhi Post Graduate Programmes IIM Bangalore IIM Ahmedabad IIM Calcutta listed among top 50 QS Masters Management Ranking 2018 released recently HEC Paris crowned best institution world studying master ''s degree management first-ever Quacquarelli Symonds QS World University Rankings Masters Management Rankings New 2018 management ranking one three QS Business Masters Rankings rankings also published finance business analytics.Indian Institute Technology Bombay IIT Bombay secured top position India moved rank 9 prestigious QS University Rankings BRICS declared last week.Indian Institute Management IIM Bangalore placed 22nd rank IIM Ahmedabad secured next 23rd position IIM Calcutta ''s PGM ranked 46.IIM Calcutta Indian member CEMS alliance leading global business schools multinational corporations placed 8th ranking CEMS Master ''s International Management programme googletag.cmd.push function googletag.display `` adslotNativeVideo In ''Global MBA course IIM Ahmedabad 49 could secure position top 50 IIM Bangalore Indian School Business top 100 58th 93th positions Regionally however composed Asia Australia New Zealand three colleges ranked among top 20 seventh ninth 12th positions IIM Calcutta top 100 also found place 18 regionally Global MBA course.One five French schools top 20 management HEC Paris ranked ahead London Business School ESADE second third respectively `` 29 different countries represented inaugural ranking excluding various global locations CEMS ranked eighth world The United Kingdom most-represented country 24 schools closely followed 21 United States said QS rankings.In rankings universities schools assessed measuring five different metrics slight differences weighting ranking Those metrics Employability Alumni Outcomes Value Money Thought Leadership Class Faculty Diversity.Union Human Resources Development minister Prakash Javadekar congratulated Indian institutions achievement Congratulations IIM_Bangalore 22 IIMAhmedabad 23 IIM_Calcutta 46 ranked among top 50 Masters Management courses worlduniranking 2018 QSWUR- Prakash Javadekar PrakashJavdekar November 29 2017 QS Ranking 2018 Top 10 Masters Management ProgrammesHEC Paris ranked best world1 HEC Paris ''s Strategic Management Programme France 2 London Business School ''s Management Programme United Kingdom 3 ESADE International Management Programme Spain 4 ESSEC ''s Strategy Management International Business Programme France 5 Imperial College London ''s Management Programme United Kingdom 6 IE Business School ''s Management Programme Spain 7 London School Economics ''s Management Strategy Programme United Kingdom 8 CEMS ''s International Management Programme9 University Copenhagen ''s International Management Programme Denmark 10 University St. Gallen ''s Business Management Programme Switzerland With Inputs IANS CommentsClick Education News'
$ws.Range('C6').Value = 'news'
$ws.Range('B7').Value = 'As per report International Data Corporation IDC Xiaomi held maximum market share top 50 Indian cities Q3 2017 period followed Samsung Lenovo including Moto Oppo Vivo The top 50 Indian cities accounted 50 percent overall smartphone sales country third quarter year Among major cities Tier cities country considered `` key volume driver Indian smartphone market sequential quarterly growth 29 percent last quarter IDC also reveals Bhopal Gurugram Jaipur Tier II cities emerged `` fastest growing cities reporting massive growth 40 percent previous quarter.IDC India Monthly City-Level Smartphone Tracker claims top 50 cities Xiaomi ousted Samsung take 26.5 percent market It said registered 120 percent growth period Redmi Note 4 contributed 40 percent volumes becoming highest selling smartphones across cities Offline expansion company ''s Preferred Partner programme Mi Stores added growth IDC said.Samsung moved second slot 24.1 percent market share saw 15 percent quarterly growth IDC says brand continues dominate several cities thanks `` deep penetration reach The Samsung Galaxy J2 Galaxy J7 Max Galaxy J7 Nxt made 50 percent volumes South Korean manufacturer IDC said.Lenovo including Moto brand third spot 10.3 percent share grew 8 percent previous quarter New Delhi Mumbai Bengaluru top cities company Lenovo K6 Power Moto E4 Plus Lenovo K6 Note top sellers.Oppo moved fourth position 5.2 percent share Oppo F3 F3 Plus A71 smartphones major contributors growth company enjoyed Finally Vivo slipped fifth spot 5.1 percent said seen marginal decline sales last quarter IDC cites `` reduced channel spends primary reason fall optimistic success Vivo V7+ Vivo Y69 smartphones bring back fourth position.IDC says metro cities continued major contributors growing smartphone sales India Alone Delhi Mumbai contributed one-fourth total sales top 50 cities However expansion data supply 4G penetration across country enabled Tier Tier II cities attract smartphone vendors `` Almost vendors decided go aggressive offline channel expansion It interesting see able match distribution strength retail spends offline dominant players like Samsung Oppo Vivo Considering huge hitherto tapped demand smaller cities towns eventually vendors forced tweak strategies reach customers markets said Navkendar Singh associate research director IDC India press statement.Vendors reportedly diversifying sales marketing strategies e-tailers like Amazon Flipkart expanding presence India cater growing demand smartphones Furthermore IDC claims e-commerce companies contributing around 40 percent market top 50 cities country `` Mega-online sale prominent e-trailers acted catalyst robust growth across city tiers e-tailers contributing around 40 percent market top 50 cities India However offline channel saw average traction festive season said Himanshu Jain market analyst IDC India.A previous IDC report revealed smartphone shipments India reached record 39 million units last quarter 40 percent previous quarter 21 percent period last year Xiaomi surfaced top smartphone brand country jointly sharing leadership Samsung market share 23.5 percent third quarter. div class=''compare-product-widget /div'
$ws.Range('B8').Value = 'Researchers Princeton ''s Center Information Technology Policy CITP claim 400 world ''s top 50,000 websites use ''session replay scripts track user behaviour While may disconcerting researchers add sites often strip personally identifiable user information behaviour data glean potentially giving hackers access trove personal data sometimes even including passwords data exposed.Detailing findings last week first several posts online privacy CITP researchers Steve Englehart Gunes Acar Arvind Narayan said looked seven top session replay companies provide session replay scripts frameworks websites These namely Clicktale FullStory Hotjar SessionCam Smartlook UserReplay Yandex To scrutinise data collected collection took place researchers set test pages session replay scripts six above-mentioned companies They also able estimate number popular sites use scripts.The researchers claim least 482 world ''s top 50,000 websites use session replay scripts number may lower side scripts n''t record actions every user visits throwing researchers detection rate Researchers compiled full list script-using websites found Getting bit business practice backfire users researchers say host information usually ends collected session linked personally identifiable data `` Collection page content third-party replay scripts may cause sensitive information medical conditions credit card details personal information displayed page leak third-party part recording This may expose users identity theft online scams unwanted behavior The true collection user inputs checkout registration processes CITP researchers explain.Some session replay script providers like SessionCam UserReplay n''t collect user data instead tracking clicks almost provide dashboard automatic manual redaction tools remove user data However remain problems approach user data still usually ends collected due sheer volume making manual scrubbing infeasible content displayed screen always collected This last especially worrying oftentimes even sites user data redaction methods place end collecting displayed content case Walgreens contained user names medical conditions prescriptions.Finally websites hosting session replay scripts may protected encrypted HTTPS protocol session replay dashboards may use vulnerable HTTP like provided Hotjar Smartlook Yandex CITP researchers noted HTTP would allow attackers use man-in-the-middle attacks get access user data transmitted third-party servers Yandex statement Motherboard responded claims said `` HTTP used intentionally session recordings load websites using iframe Unfortunately loading HTTP content HTTPS websites prohibited browser level HTTP player required support HTTP websites feature `` Among sites use session replay scripts major names include Bonobos Fidelity apart already named Walgreens After publication CITP study last week Bonobos told Wired ended data sharing FullStory reviewing protocols better protect user data A Fidelity spokesperson told Motherboard protection customer data highest priority n''t clarify would stop using scripts Walgreens took tack Bonobos said `` abundance caution stopped sharing data FullStory investigated claims.The study notes ad-blocking lists tracking protection services like EasyList EasyPrivacy provide measure safety block everything Motherboard reports Adblock Plus updated post publication CITP study block named scripts. div class=''compare-product-widget /div'
$ws.Range('C8').Value = 'gadgets'
$ws.Range('B9').Value = 'The Ambani family led brothers Mukesh Ambani Anil Ambani claimed top slot Forbes list 50 richest Asian families net worth 44.8 billion 19 billion previous year India enjoys biggest presence list 18 families among region ''s top 50 The Ambanis superseded Lees Samsung empire claim No 1 spot The wealth oil-to-telecom tycoon Mukesh Ambani also richest Indian surged rally shares Reliance Industries strong customer response Reliance Jio His younger brother Anil Ambani heads Reliance Communications.Here 10 things must know Forbes list Asia ''s richest families:1 The combined wealth Ambanis soared whopping 19 billion pipped South Korea ''s Lee family claim top notch rank Forbes list Asia ''s wealthiest clans.2 `` Because India ''s size Ambanis never dominant Samsung ''s Lees Korea But Mukesh ''s next generation establishing presence Reliance Jio telecom operation story could play years Tim Ferguson Editor Forbes Asia said Forbes website forbesindia.com.3 The Ambanis Indian family top 10 richest families list Asia `` No family highlights surge better Ambanis India year ''s biggest gainer dollar percentage terms Forbes said googletag.cmd.push function googletag.display `` adslotNativeVideo 4 Other richest Indian families list feature Premjis rank 11 19.2 billion net worth Hindujas 12th 18.8 billion Mittals 14th 17.2 billion Mistrys 16th 16.1 billion Birlas 19th 14.1 billion .5 Other Indian families got exalted list include Godrej family 20th 14 billion Bajajs 26th 9.3 billion Jindal family 32nd 7.68 billion Dabur ''s Burmans 35th 7.05 billion Eicher Motors Lals 36th 7 billion Shree Cement ''s Bangur family 37th 6.7 billion .6 Motherson Sumi Systems Sehgal family 41st 6.15 billion Wadia family 42nd 6.14 billion DLF ''s Kushal Pal Singh 44th 6.1 billion Patel family controls Cadila 45th 6 billion Piramals 47th 5.38 billion Munjals 48th 5.37 billion also made list.7 South Korea ''s Lee family ranked second Forbes list saw wealth soar 11.2 billion 40.8 billion year shares Samsung Electronics moved 75 per cent past year.8 Hong Kong ''s Kwok family controls Sun Hung Kai Properties ranked third year net wealth 40.4 billion Comments9 The list Asia ''s 50 Richest Families snapshot wealth using stock prices currency exchange rates close markets November 3 Private companies valued using financial ratios comparisons similar publicly traded firms.10 The ticket entry year ''s list 5 billion 1.6 billion 2016'
$ws.Range('C9').Value = 'news'
$ws.Range('B10').Value = 'Sania Mirza dropped top 10 first time since March 2014 © AFP Highlight Story widget starts Highlight Story widget end A consistent Divij Sharan broke top 50 career-high ranking Sania Mirza dropped top 10 first time since March 2014 per latest ATP WTA charts released Monday The left-handed Sharan reached final last three tournaments participated winning ATP 250 European Open ending runner-up Challenger events Brest Tashkent different partners He ranked exactly 50th gaining place The 31-year-old Delhi player let split long-time partner Purav Raja affect consistently well Tour despite stable partnership googletag.cmd.push function googletag.display ''adslotNativeVideo Now India two top 50 players doubles Rohan Bopanna 15 leading chart.Purav Raja 62 Leander Paes 70 Jeevan Nedunchezhiyan 97 top-100 Indian players `` made first ATP point 2003 It ''s long exciting journey since There lot hard-work sacrifices experiences travel decisions ups downs importantly grateful opportunity love playing tennis Sharan told PTI Bratislava Slovakia Sharan 19 ITF Futures titles 13 Challenger trophies three ATP Tour titles said `` belief win Grand Slam `` The guys losing close matches ones winning Slams We beaten sometimes feel win Grand Slam said `` employed Indian Oil since 2004 support massive career Head Adidas also providing racquets kit last four years added Sharan said next goal make top 30 play ATP 1000s Master series events `` It ''s good motivation encouragement Sharan reached pre-quarterfinals French Open year said.Also making good progress doubles Vishnu Vardhan jumped 16 places touch career-best ranking 118 Vardhan partner N Sriram Balaji 139 Shenzhen Challenger Saturday earn 90 ranking points each.In singles Yuki Bhambri 140 continues lead chart followed Ramkumar Ramanathan 148 Prajnesh Gunneswaran 255 Sumit Nagal 331 Balaji 350 .In WTA charts Ankita Raina India ''s top-ranked player number 281 followed Karman Kaur Thandi 307 Pranjala Yadlapalli 483 In doubles Sania placed number 12 following loss three places.Sania figured consistently top-10 since March 17 2013 The last time featured outside top-10 March 3 2013.The Hyderabadi world number one doubles entire 2016 season Her slide started season began number one could qualify year-end WTA Finale top-eight pairs top-eight singles players compete Comments Topics Sania Mirza Divij Sharan Tennis Get latest FIFA World Cup 2018 news check World Cup 2018 schedule Football live score World Cup standings Like us Facebook follow us Twitter FIFA World Cup 2018 updates You also download NDTV Cricket app Android iOS'
$ws.Range('C10').Value = 'sports'
$ws.Range('B11').Value = 'Mumbai The government working World Bank recognise 200 reforms help propel India top-50 bracket ease business according Department Industrial Policy Promotion DIPP secretary Ramesh Abhishek `` We already implemented 122 reforms year working World Bank recognise We initiate 90 ease business reforms year Mr Abhishek told reporters sidelines CII Invest North conclave.Yesterday India jumped 30 places rank 100th World Bank ''s ''ease business ranking helped slew reforms taxation licensing investor protection bankruptcy resolution `` A jump 30 positions World Bank rankings extremely commendable Now aim move top-50 countries bracket added.Mr Abhishek said department already started meeting stakeholders taking feedback reform measures taken government improve business climate `` This exercise helped us lot This time focused mainly taking feedback stakeholders All nodal ministries taking feedback said.googletag.cmd.push function googletag.display `` adslotNativeVideo The Secretary also said World Bank report acknowledged GST Goods Services Tax substantial reform `` hopefully positive impact country ''s ranking next year CommentsHe urged states address issues like inadequate infrastructure inflexible labour policies make competitive `` States invest skilling need come policies fixed term employment like government worked apparel sector We also working similar package leather footwear sector added Mr Abhishek'
$ws.Range('B12').Value = 'New Delhi Buoyed jump India ''s ease business ranking Finance Minister Arun Jaitley today vowed continue reforms help break top 50 nations.According World Bank India ''s rank ''ease business scale risen 130th 100th year helped slew reforms taxation licensing investor protection bankruptcy resolution.Addressing press conference minutes World Bank ranking made public Mr Jaitley said India major country named pursuing structural reforms `` In 2014 142nd last two years improved 131st 130th These generalised rankings It happened specific areas take tough parameters ranking said.The `` highest jump ranking possible significant improvements 10 judging parameters made last 3-4 years `` becomes easy business India said.The World Bank said named 10 countries structural reforms taken `` India one `` googletag.cmd.push function googletag.display `` adslotNativeVideo Mr Jaitley said 3 years ranking improved 142nd 100th sufficient work ongoing areas India lags behind `` There reason believe capacity significantly improve position said.The World Bank says India still lags areas starting business enforcing contracts dealing construction permits.It takes 30 days register new business 127 days 15 years ago `` number procedures still cumbersome local entrepreneurs still need go 12 procedures said `` The Prime Minister said target take India immediately first 50 believe doable therefore remaining 3-4 areas work done pushing greater force said.Elaborating ranking said significant improvement parameters protecting minority investors interest changes made company law `` We ranked 4th world This biggest jump field `` In availability credit India ranked 29th `` Biggest jump improvements done taxation front Till last year 189 nations 172nd It means paying taxes cumbersome process But year moved 53 places said.Also enactment new bankruptcy law India ''s ranking insolvency resolution also improved.The ranking include impact implementation Goods Services Tax GST created cut judging business environment June 1 `` In coming years expecting get credit said.Besides areas government done lot work work still progress get reflected numbers coming 1-2 years said.Giving examples said India ranked low 164th enforcement contracts arbitration act help improve this.Similarly states impressed upon give building construction permits online fast permissions handed `` One reform factored starting new business In 156th position created provision company act ... reform fully factored ranking sure future said.India ranked 156th latest ranking `` Work progress think able expedite work next months areas huge scope us jump 100 said.The finance minister said government seriously taking 3-4 issues India still lags behind next 1-2 years improvement 20-30 points ranking `` If see first available data August September think every data coming seems positive zone said CommentsMr Jaitley said reforms continuous process `` The 3-4 areas n''t treated high areas significant amount work progress impact shown ground ''m sure taken consideration next year ''s ranking said'
$ws.Range('B13').Value = 'Here ten points story Sharing details report Finance Minister Arun Jaitley said Tuesday India made biggest jump index `` The Prime Minister said top 50 The fact three years jumped 142 100 believe doable Mr Jaitley said.PM Modi series tweets said jump India ''s ranking result all-round multi-sectoral reform push `` Easier business environment leading historic opportunities entrepreneurs particularly MSME medium small scale sector sector bringing prosperity tweeted.The report based data New Delhi Mumbai ranked India among top 10 `` improvers globally done better eight 10 business indicators Mr Jaitley said government trying improve 10.India jumped 53 places 172 119 ease paying taxes leapt 4th rank protecting minority investors moved 15 places 44 29 ease businesses getting credit.There also red flags India fallen 16 places ease registering property slid slightly trading across borders one rank 155 156 ease starting businesses Mr Jaitley said would improve next year since report n''t factored changes made law Companies Act.New Zealand tops World Bank index ease business year Singapore number 2 Denmark placing third Somalia last China 78 Pakistan 147 Bangladesh 177 The report boost government defend economic policies amid criticism GDP growth rate hitting three-year low April June quarter.The government insisted fundamentals economy strong rejected allegations two big reforms last year ''s demonetisation new national tax GST launched year blame economic slowdown saying deliver big benefits `` Today ''s result clear signal India rest world country ready open business many decades competing preferred place business globally Annette Dixon World Bank ''s vice president South Asia told reporters New Delhi.The opposition Congress however kept attack saying NDA government ''s decisions like `` cease businesses medium small sectors The hype built government around report Congress said attempt seek `` manufactured relevance World Bank people suffering due GST notes ban.googletag.cmd.push function googletag.display `` adslotNativeVideo Comments'
$ws.Range('B14').Value = 'Mumbai October 17 Indian Institute Technology Bombay IIT Bombay moved one step 35th rank last year 34th rank year Quacquarelli Symonds QS Asia University Rankings Among top 50 Overall Institution improved 2016/17 performance one place Indian Institute Technology IIT Delhi ranked 41 IIT Madras placed rank 48 QS Asia University Rankings 2018 IISc Bangalore IIT Kanpur IIT Kharagpur University Delhi IIT Roorkee IIT Guwahati also found place 100 year.Earlier June year IIT Bombay ranked 179th QS World University Rankings showing improvement 2017 performance 40 places.The results Asia University Rankings release Monday October 16 2017 QS British company IIT Bombay ranks one India overall score 75.7 100 The Institute score 88.1 academic reputation 97.5 employer reputation 44.5 faculty student ration 96.3 staff PhD 89.6 papers per faculty 9.7 international faculty 6.3 international students 10.9 inbound exchange students 9.1 outbound exchange students 41.4 citations per paper Among ten parameters Employer Reputation indicated strongest one IIT Bombay rank 21 globally `` glad IIT Bombay continues rise rankings This reflection excellent work done faculty students Institute Director Instituet Prof. Devang Khakhar said.googletag.cmd.push function googletag.display `` adslotNativeVideo Considering approximately 11,900 universities globally made IIT Bombay one top 1.0 universities Asia.In QS Ranking 2018 released June year IIT Delhi dethroned IISc Bengaluru In 2017 rankings IISC placed highest among Indian institute managed rank top 200 IISc featured 152nd rank 2017 rankings QS Asia University Rankings 2018 Top 10Nanyang Technological University Singapore NTU Singapore Asia ''s top university according QS Asia University NTU renowned among world ''s best universities engineering technology rose two spots last year placed leading position QS Asia University league table ahead 400 universities Asia QS Asia University Rankings 2018 NTU Singapore ranked top list NTU made fantastic progress risen No 14 QS Asian rankings began 2009 The latest rankings outstanding achievement strong endorsement NTU ''s excellent reputation built years wish thank entire NTU community contributions NTU President Professor Bertil Anderson said.Universities Singapore Hong Kong China dominated rankings:1 Nanyang Technological University NTU Singapore 2 National University Singapore NUS 3 The Hong Kong University Science Technology Hong Kong 4 KAIST Korea Advanced Institute Science Technology South Korea 5 The University Hong Kong Hong Kong 6 Tsinghua University China 7 Fudan University China 8 City University Hong Kong Hong Kong 9 Peking University China 10 The Chinese University Hong Kong CUHK Click Education News Comments'
$ws.Range('B15').Value = 'Mumbai Banks may write 60 percent value bad loans 50 large stressed asset accounts take hit Rs 2.4 lakh crore ratings agency CRISIL said Wednesday.The 50 stressed companies account Rs 4 lakh crore soured loans largely metals construction power industries account nearly half total non-performing assets banking sector end March CRISIL said.The credit agency part S P Global estimated banks provisioned around 40 percent exposure assets CommentsBanks total non-performing loans Rs 7.29 lakh crore 5 percent India ''s gross domestic product end-March.The report comes rival India Ratings Research estimated banks taking 12 country ''s largest defaulters bankruptcy court would need make additional provisioning least Rs 18,000 crore.© Thomson Reuters 2017googletag.cmd.push function googletag.display `` adslotNativeVideo'
$ws.Range('B16').Value = 'New Delhi Enthused success SBI merger Finance Ministry considering clearing another proposal public sector banking space fiscal end goal create 4-5 global sized lenders.Five associate banks Bharatiya Mahila Bank became part SBI April 1 2017 catapulting country ''s largest lender among top 50 banks world.Now Finance Ministry looking replicate model case state-run banks reach critical mass compete global peers `` Consolidation must ... decision regard would based commercially prudent parameters If NPA situation gets better could one merger towards end fiscal senior official told PTI.Toxic loans public sector banks rose Rs 1 lakh crore Rs 6.06 lakh crore April-December 2016-17 bulk came power steel road infrastructure textile sectors.Finance Minister Arun Jaitley several occasions said India needs 5-6 banks global size scale consolidation banking sector done appropriate time.googletag.cmd.push function googletag.display `` adslotNativeVideo Whenever consolidation happens takes consideration interest stakeholders including employees shareholders official said adding win-win parties.A balancing act done merger given clearance various authorities regulators official added.The future merger proposals banking sector also require clearance Competition Commission India CCI see merged entity going create monopoly.In last consolidation drive saw light day April CCI nod needed case merger Bharatiya Mahila Bank BMB SBI There requirement merger associate banks SBI part parent.The Finance Ministry sought help government think-tank NITI Aayog global consultancy firms examine possibility next round consolidation PSU banks aim create lenders global size scale.NITI Aayog ''s report expected set tone roadmap consolidation future.There factors like regional balance geographical reach financial burden smooth human resource transition looked taking merger decision official said.The official added merger weak bank strong bank `` could pull latter `` There low-hanging fruits Big lenders like Bank Baroda take turnaround banks southern region Indian Overseas Bank Dena Bank could merged large South Indian bank official explained.Five associates BMB became part SBI April 1 2017 catapulting country ''s largest lender among top 50 banks world.State Bank Bikaner Jaipur SBBJ State Bank Hyderabad SBH State Bank Mysore SBM State Bank Patiala SBP State Bank Travancore SBT besides BMB merged SBI.With merger total customer base SBI reached around 37 crore branch network around 24,000 nearly 59,000 ATMs across country The merged entity began operation deposit base Rs 26 lakh crore advances Rs 18.50 lakh crore CommentsThe government February approved merger five associate banks SBI Later March cabinet approved merger BMB well SBI first merged State Bank Saurashtra 2008 Two years later State Bank Indore merged This story edited NDTV staff auto-generated syndicated feed'
$ws.Range('B17').Value = 'Jaipur Rajasthan performed spectacularly NEET National Eligibility cum Entrance Tests exam year Out 79,057 appeared NEET 45,336 candidates cleared test 182 Jaipur alone four top 50 state.Prince Chaudhary small village Barmer The 17-year-old became famous overnight securing fifth rank country Son medical store owner Prince wants pursue medical course All India Institute Medical Science Delhi become neuro physician `` happy shy teenager told NDTV And though studied Hindi medium school till class 10 two years intense coaching career institute Kota gave confidence take exam English scoring 686 720.Not Prince one every five students polished exam skills coaching centres Kota claims one teachers institute `` Coaching institutes create awareness able help students focused preparation said Naveen Maheshwari heads one sought coaching institutes Kota The city famous cotton sarees seen mushrooming coaching centres thriving NEET boom There around 1.5 lakh students 150 centres city.googletag.cmd.push function googletag.display `` adslotNativeVideo Students mostly enroll institutes class 11 spend two years studying exactly crack exam says Manan Bhardwaj Jaipur got India rank 150 `` It matter much study matters study Manan told NDTV Comments We seeing child Barmer cleared exam also got India rank 5 unheard till years ago Last year saw Rupa Yadav child bride married 8 clear exam So certainly awareness coaching help children back ground Mr Maheshwari said.In coaching capital 150 institutes annual turnover Rs 1500 crore say reports And ''s tuition fees lakhs aspirants come Uttar Pradesh Madhya Pradesh Chhattisgarh Bihar parts northern India rent rooms duration courses There booming PG accommodation business private houses dormitories facilities'
$ws.Range('C17').Value = 'news'
$ws.Range('B18').Value = 'Netflix made name back really good TV shows continues spend billions dollars year build portfolio original series prominently showcases homepage attract eyeballs self-produced efforts By tends background popular critically-acclaimed series n''t make arguably better anything Netflix made history.There ''s need highlight casual viewers may miss Netflix ''s ''ve taken upon gather best TV shows streaming Netflix To prepare list used aggregate ratings Rotten Tomatoes Metacritic IMDb draw shortlist picked top 50 Here best TV shows Netflix India sorted alphabetically This list updated every month worthy additions movies removed service bookmark page keep checking in.The 100 Best Movies Netflix India A Series Unfortunate Events 2017 – Present Three resourceful orphans – inventor reader baby sharp teeth – must outsmart evil distant relative ''s sizeable fortune trying uncover mystery behind secret society parents involved Based Lemony Snicket ''s popular series children ''s books Alias Grace 2017 Margaret Atwood ''s 1996 novel name 19th-century Canadian woman convicted double murder becomes subject criminal psychologist – profession n''t exist name – adapted screen six-part miniseries American Crime Story 2016 – Present A true crime anthology series prolific producer Ryan Murphy follows well-known events dominated US media trial former sportsman O.J Simpson assassination fashion designer Gianni Versace Better Call Saul 2015 – Present This spin-off prequel Breaking Bad follows small-time lawyer Bob Odenkirk tendencies con artist transforms morally-challenged criminal lawyer knew Saul Goodman Black Mirror 2011 – Present Charlie Brooker ''s anthology series consisting standalone episodes – means ever-changing cast addition new settings storylines – explores unanticipated consequences new technologies often dark satirical ways BoJack Horseman 2014 – Present Set world humans anthropomorphic animals live alongside washed-up sitcom star plans comeback fame help ghostwriter ex-girlfriend ''s also agent freeloading roommate dealing rival ''s dating ghostwriter Breaking Bad 2008 – 2013 Diagnosed lung cancer struggling high school chemistry teacher Bryan Cranston decides get business making selling meth secure family ''s financial future help former student Aaron Paul Broadchurch 2013 – 2017 A look violent crimes affect small seaside town Britain eyes two investigators David Tennant Olivia Colman impact media attention suspicion grief close-knit community Brooklyn Nine-Nine 2013 – Present The lives group detectives fictional New York precinct – ensemble featuring Any Samberg AndreBraugher Terry Crews – get sitcom treatment The Office co-creator Michael Schur Crazy Ex-Girlfriend 2015 – Present A successful young female lawyer Rachel Bloom also co-creator suffers depression anxiety uproots life New York moves suburb California find love happiness The Crown 2016 – Present A look life Britain ''s longest-ruling monarch Queen Elizabeth II wedding 1947 present day including political rivalries romances world-changing events reign Daredevil 2015 – Present A blind man fights injustice attorney day masked vigilante night present-day New York struggling deal Catholic guilt arises actions Dear White People 2017 – Present Set fictional predominantly white Ivy League college Justin Simien turns successful satirical film broader series black students dealing social injustice figuring The End F***ing World 2017 – Present Two teenage outsiders – budding psychopath rebel hungry adventure – embark road trip search latter ''s real father become involved series increasingly violent events The Expanse 2015 – Present Hundreds years future mankind colonised Solar System brink war ''s crew different origins – Earth Mars Asteroid Belt – expose greatest conspiracy Fargo 2014 – Present The Coen brothers acclaimed film name inspiration black comedy/crime anthology series featuring quirky characters across different eras dealing deception intrigue murder amidst cold American Midwest Fauda 2015 – Present A former expert Israeli agent comes retirement hunt Palestinian militant show cleverly following sides conflict drawing real world events Fullmetal Alchemist Brotherhood 2009 – 2010 This direct adaptation famous manga centres two brothers looking Philosopher ''s Stone restore bodies failed experimental attempt resurrect mother goes awry But ''re ones stone GLOW 2017 – Present An exploration real-life 1980s women ''s professional wrestling – Gorgeous Ladies Of Wresting GLOW – focusing personal professional lives fictionalised Hollywood misfits involved Godless 2017 In 19th century murderous outlaw notorious gang hunt ex-protégé betrayed find conflict town gives refuge whose residents mainly women A seven-episode miniseries The Good Place 2016 – Present The Office co-creator Michael Schur spun sitcom web afterlife series following woman Kristen Bell ''s mistakenly assigned Heaven-like utopia tries better person hide The Good Wife 2009 – 2016 After humiliating sex corruption scandal puts husband behind bars wife – former state ''s attorney – must return work provide family battling unwanted spotlight Hannibal 2013 – 2015 Cancelled three seasons cult horror hit explores relationship forensic psychiatrist Mads Mikkelsen patient young FBI criminal profiler uncanny ability empathise serial killers Happy Valley 2014 – Present A British police sergeant dealing personal problems investigates cases small West Yorkshire town Two seasons six episodes possibly way House Cards 2013 – Present A betrayed US Congressman works equally conniving wife climb political ladders Washington anything needed – manipulation deceit even murder – achieve goals The Inbetweeners 2008 – 2010 A coming-of-age sitcom whose success led two movies following four British teenagers final year school series misadventures involving uncaring school staff male bonding failed sexual encounters Jane Virgin 2014 – Present Gina Rodriguez Annihilation stars devout Catholic working young Latina virgin becomes pregnant accidental artificial insemination rom-com satire end fifth season 2019 Jessica Jones 2015 – Present Suffering PTSD super-powered woman rebuilds life private investigator New York forced battle past demons every step way wants look past Lost 2004 – 2010 The survivors plane crash must work together survive mysterious tropical island somewhere South Pacific Ocean battling supernatural sci-fi elements Faltered midway recovered ended disappointing fashion Luther 2010 – Present Idris Elba stars dedicated brilliant British detective tries keep grip personal life dealing psychological fallouts crimes ''s tasked solve Mad Men 2007 – 2015 Set 1960s New York slow-burn drama offers peek inside fictional ad agency focusing one extremely talented executives Jon Hamm ''s bored simple personal life Master None 2015 – Present Loosely based Aziz Ansari ''s life experiences 30-year-old struggling actor tries get life order personal professional fronts affected Indian ancestry Mindhunter 2017 – Present In late 1970s two FBI agents push superiors expand research criminal science involves getting close personal imprisoned serial killers understand think Monty Python ''s Flying Circus 1969 – 1974 The influential British comedy group got start surreal sketch series targeted life island intellectual fashion humour wide-ranging unique gave birth term `` Pythonesque Narcos 2015 – Present A gripping look violent powerful drug cartels Colombia including infamous Pablo Escobar thecorroborative efforts various law enforcement whose job bring Orange Is New Black 2013 – Present The lives incarcerated women minimum-security federal prison upstate New York including normally law-abiding privileged woman ''s sentenced decade-old crime Orphan Black 2013 – 2017 A con artist Tatiana Maslany assumes identity woman committed suicide looked like pulled conspiracy learns ''s clone Outlander 2014 – Present Diana Gabaldon ''s best-selling books married nurse World War II ''s transported back time 1743 finds caught Jacobite risings two different men Peaky Blinders 2013 – Present The exploits Shelby crime family Birmingham England two World Wars elements borrowed 19th-century gang name legend goes used sew razor blades caps Peep Show 2003 – 2015 The lives two different dysfunctional twenty-somethings sharing flat South London day-to-day farcical antics made cult favourite Rick Morty 2013 – Present An adult animated sci-fi series cynical alcoholic mad scientist Rick weak-willed easily influenced grandson Morty go series misadventures across dimensions universes Sherlock 2010 – Present Benedict Cumberbatch Martin Freeman play famous detective doctor sidekick modern-day adaptation Sir Arthur Conan Doyle ''s stories Terrific compelling first years though derailed fourth season Star Trek 1966 – 1969 Gene Rodenberry ''s original series adventures starship 2260s led Kirk Spock McCoy n''t age well visually storytelling endures set aside major dip quality third season Star Trek The Next Generation 1987 – 1994 Rodenberry took vision even far-off future – 24th century – focusing new generation Starfleet officers new ship recognisable name Stranger Things 2016 – Present A nostalgia-heavy love letter 80s set suburban small town secret scientific experiments paranormal supernatural times human subjects starts wreak havoc normal lives residents That ''70s Show 1997 – 2006 The comedy fueled Hollywood careers stars Mila Kunis Ashton Kutcher group six high school friends going usual teenage problems trying figure identities Trollhunters 2016 – Present Guillermo del Toro turns love monster stories animated series ages following teenage boy stumbles onto mystic amulet must protect world humans trolls Unbreakable Kimmy Schmidt 2015 – Present Rescued doomsday cult 15 years young woman armed positive attitude decides start new life New York help gay wannabe-Broadway actor roommate street-wise landlady depressed out-of-touch socialite The Walking Dead 2010 – Present Based popular comic series horror drama set post-apocalyptic future survivors search safe world overrun zombies Hit peak fifth season never recovered Wentworth 2013 – Present Locked awaiting trial attempted murder husband woman adjusts life Australian prison rises ranks Gripping till season 4 post lead actress left div class=''compare-product-widget /div'
$ws.Range('C18').Value = 'gadgets'
$ws.Range('B19').Value = 'New Delhi TIMES Higher Education THE World University Ranking released rankings educational institutes Asia two Indian institutes made way top 50 Asian institutes The two institutes made way top 50 Indian Institute Science 29th place Indian Institute Technology Bombay 44th rank Other institutes made way top 100 Indian Institute Technology Kharagpur Indian Institute Technology Roorkee Indian Institute Technology Kanpur Indian Institute Technology Delhi The ranks released THE World University Ranking Indian Institute Science IISc world rank 29Indian Institute Technology IIT Bombay world rank 44Indian Institute Technology IIT Kharagpur world rank 60Indian Institute Technology IIT Roorkee world rank 65Indian Institute Technology IIT Kanpur world rank 81Indian Institute Technology IIT Delhi world rank 86Among made top 200 Tezpur University 100th place Indian Institute Technology Madras 103rd place Indian Institute Technology Guwahati 112th place Panjab University 114th place National Institute Technology Rourkela 126th place Jadavpur University 127th place Indian School Mines 141st place University Delhi 144th place Aligarh Muslim University 158th place Savitribai Phule Pune University 188th place IIT-BHU 194th place CommentsA total 17 institutes India made way top 200 Asian Universities Though representation India increased list several universities fallen rankings comparison last year Both IISc IIT Bombay dropped two places 29th 44th position respectively IIT Madras witnessed biggest decline year fallen last year ''s 41st place 103rd pace year googletag.cmd.push function googletag.display `` adslotNativeVideo Click Education News'
$ws.Range('B20').Value = 'New Delhi A differently abled Indian teacher little-known institution Uttarakhand beat tough competition around world named today among top 50 contenders USD 1 million global award Pradeep Negi economics social science computer science teacher Government Inter College BHEL Ranipur Haridwar shortlisted annual Varkey Foundation Global Teacher Prize 2018 30,000 nominations applications 173 countries around world `` It unbelievable proud chosen one top 50 teachers It motivated teachers community said Mr Negi `` handicapped n''t stopped succeeding glad inspire disabled people India developed lot innovative work ICT poor children glad recognised said Mr Negi physically disabled age two due polio.Read Maggie McDonnell Teacher From Canada Wins Global Teacher Prize 2017Negi chosen fighting spirit suffered bullying school disability overcoming adversities use new technology tools make economics social subjects compelling students.googletag.cmd.push function googletag.display `` adslotNativeVideo Congratulations Pradeep Negi reaching final 50 hope story inspires looking enter teaching profession shines spotlight incredible work teachers world every day said Sunny Varkey founder Varkey Foundation Global Teacher Prize `` We intend keep momentum going journey continues return teachers rightful position one respected professions society said.Mr Negi helped train 1,200 teachers state.The Global Teacher Prize set recognise one exceptional teacher every year made outstanding contribution profession.The top 50 shortlisted teachers narrowed 10 finalists committee results announced February 2018.The winner chosen final 10 finalists Global Teacher Prize Academy.All 10 finalists invited Dubai Award ceremony Global Education Skills Forum GESF March 18 2018 winner announced live stage CommentsClick Education News Except headline story edited NDTV staff published syndicated feed'
$ws.Range('B21').Value = 'New Delhi Post Graduate Programmes IIM Bangalore IIM Ahmedabad IIM Calcutta listed among top 50 QS Masters Management Ranking 2018 released recently HEC Paris crowned best institution world studying master ''s degree management first-ever Quacquarelli Symonds QS World University Rankings Masters Management Rankings New 2018 management ranking one three QS Business Masters Rankings rankings also published finance business analytics.Indian Institute Technology Bombay IIT Bombay secured top position India moved rank 9 prestigious QS University Rankings BRICS declared last week.Indian Institute Management IIM Bangalore placed 22nd rank IIM Ahmedabad secured next 23rd position IIM Calcutta ''s PGM ranked 46.IIM Calcutta Indian member CEMS alliance leading global business schools multinational corporations placed 8th ranking CEMS Master ''s International Management programme googletag.cmd.push function googletag.display `` adslotNativeVideo In ''Global MBA course IIM Ahmedabad 49 could secure position top 50 IIM Bangalore Indian School Business top 100 58th 93th positions Regionally however composed Asia Australia New Zealand three colleges ranked among top 20 seventh ninth 12th positions IIM Calcutta top 100 also found place 18 regionally Global MBA course.One five French schools top 20 management HEC Paris ranked ahead London Business School ESADE second third respectively `` 29 different countries represented inaugural ranking excluding various global locations CEMS ranked eighth world The United Kingdom most-represented country 24 schools closely followed 21 United States said QS rankings.In rankings universities schools assessed measuring five different metrics slight differences weighting ranking Those metrics Employability Alumni Outcomes Value Money Thought Leadership Class Faculty Diversity.Union Human Resources Development minister Prakash Javadekar congratulated Indian institutions achievement Congratulations IIM_Bangalore 22 IIMAhmedabad 23 IIM_Calcutta 46 ranked among top 50 Masters Management courses worlduniranking 2018 QSWUR- Prakash Javadekar PrakashJavdekar November 29 2017 QS Ranking 2018 Top 10 Masters Management ProgrammesHEC Paris ranked best world1 HEC Paris ''s Strategic Management Programme France 2 London Business School ''s Management Programme United Kingdom 3 ESADE International Management Programme Spain 4 ESSEC ''s Strategy Management International Business Programme France 5 Imperial College London ''s Management Programme United Kingdom 6 IE Business School ''s Management Programme Spain 7 London School Economics ''s Management Strategy Programme United Kingdom 8 CEMS ''s International Management Programme9 University Copenhagen ''s International Management Programme Denmark 10 University St. Gallen ''s Business Management Programme Switzerland With Inputs IANS CommentsClick Education News'
$ws.Range('C21').Value = 'news'
$ws.Range('B22').Value = 'As per report International Data Corporation IDC Xiaomi held maximum market share top 50 Indian cities Q3 2017 period followed Samsung Lenovo including Moto Oppo Vivo The top 50 Indian cities accounted 50 percent overall smartphone sales country third quarter year Among major cities Tier cities country considered `` key volume driver Indian smartphone market sequential quarterly growth 29 percent last quarter IDC also reveals Bhopal Gurugram Jaipur Tier II cities emerged `` fastest growing cities reporting massive growth 40 percent previous quarter.IDC India Monthly City-Level Smartphone Tracker claims top 50 cities Xiaomi ousted Samsung take 26.5 percent market It said registered 120 percent growth period Redmi Note 4 contributed 40 percent volumes becoming highest selling smartphones across cities Offline expansion company ''s Preferred Partner programme Mi Stores added growth IDC said.Samsung moved second slot 24.1 percent market share saw 15 percent quarterly growth IDC says brand continues dominate several cities thanks `` deep penetration reach The Samsung Galaxy J2 Galaxy J7 Max Galaxy J7 Nxt made 50 percent volumes South Korean manufacturer IDC said.Lenovo including Moto brand third spot 10.3 percent share grew 8 percent previous quarter New Delhi Mumbai Bengaluru top cities company Lenovo K6 Power Moto E4 Plus Lenovo K6 Note top sellers.Oppo moved fourth position 5.2 percent share Oppo F3 F3 Plus A71 smartphones major contributors growth company enjoyed Finally Vivo slipped fifth spot 5.1 percent said seen marginal decline sales last quarter IDC cites `` reduced channel spends primary reason fall optimistic success Vivo V7+ Vivo Y69 smartphones bring back fourth position.IDC says metro cities continued major contributors growing smartphone sales India Alone Delhi Mumbai contributed one-fourth total sales top 50 cities However expansion data supply 4G penetration across country enabled Tier Tier II cities attract smartphone vendors `` Almost vendors decided go aggressive offline channel expansion It interesting see able match distribution strength retail spends offline dominant players like Samsung Oppo Vivo Considering huge hitherto tapped demand smaller cities towns eventually vendors forced tweak strategies reach customers markets said Navkendar Singh associate research director IDC India press statement.Vendors reportedly diversifying sales marketing strategies e-tailers like Amazon Flipkart expanding presence India cater growing demand smartphones Furthermore IDC claims e-commerce companies contributing around 40 percent market top 50 cities country `` Mega-online sale prominent e-trailers acted catalyst robust growth across city tiers e-tailers contributing around 40 percent market top 50 cities India However offline channel saw average traction festive season said Himanshu Jain market analyst IDC India.A previous IDC report revealed smartphone shipments India reached record 39 million units last quarter 40 percent previous quarter 21 percent period last year Xiaomi surfaced top smartphone brand country jointly sharing leadership Samsung market share 23.5 percent third quarter. div class=''compare-product-widget /div'
$ws.Range('B23').Value = 'Researchers Princeton ''s Center Information Technology Policy CITP claim 400 world ''s top 50,000 websites use ''session replay scripts track user behaviour While may disconcerting researchers add sites often strip personally identifiable user information behaviour data glean potentially giving hackers access trove personal data sometimes even including passwords data exposed.Detailing findings last week first several posts online privacy CITP researchers Steve Englehart Gunes Acar Arvind Narayan said looked seven top session replay companies provide session replay scripts frameworks websites These namely Clicktale FullStory Hotjar SessionCam Smartlook UserReplay Yandex To scrutinise data collected collection took place researchers set test pages session replay scripts six above-mentioned companies They also able estimate number popular sites use scripts.The researchers claim least 482 world ''s top 50,000 websites use session replay scripts number may lower side scripts n''t record actions every user visits throwing researchers detection rate Researchers compiled full list script-using websites found Getting bit business practice backfire users researchers say host information usually ends collected session linked personally identifiable data `` Collection page content third-party replay scripts may cause sensitive information medical conditions credit card details personal information displayed page leak third-party part recording This may expose users identity theft online scams unwanted behavior The true collection user inputs checkout registration processes CITP researchers explain.Some session replay script providers like SessionCam UserReplay n''t collect user data instead tracking clicks almost provide dashboard automatic manual redaction tools remove user data However remain problems approach user data still usually ends collected due sheer volume making manual scrubbing infeasible content displayed screen always collected This last especially worrying oftentimes even sites user data redaction methods place end collecting displayed content case Walgreens contained user names medical conditions prescriptions.Finally websites hosting session replay scripts may protected encrypted HTTPS protocol session replay dashboards may use vulnerable HTTP like provided Hotjar Smartlook Yandex CITP researchers noted HTTP would allow attackers use man-in-the-middle attacks get access user data transmitted third-party servers Yandex statement Motherboard responded claims said `` HTTP used intentionally session recordings load websites using iframe Unfortunately loading HTTP content HTTPS websites prohibited browser level HTTP player required support HTTP websites feature `` Among sites use session replay scripts major names include Bonobos Fidelity apart already named Walgreens After publication CITP study last week Bonobos told Wired ended data sharing FullStory reviewing protocols better protect user data A Fidelity spokesperson told Motherboard protection customer data highest priority n''t clarify would stop using scripts Walgreens took tack Bonobos said `` abundance caution stopped sharing data FullStory investigated claims.The study notes ad-blocking lists tracking protection services like EasyList EasyPrivacy provide measure safety block everything Motherboard reports Adblock Plus updated post publication CITP study block named scripts. div class=''compare-product-widget /div'
$ws.Range('C23').Value = 'gadgets'
$ws.Range('B24').Value = 'The Ambani family led brothers Mukesh Ambani Anil Ambani claimed top slot Forbes list 50 richest Asian families net worth 44.8 billion 19 billion previous year India enjoys biggest presence list 18 families among region ''s top 50 The Ambanis superseded Lees Samsung empire claim No 1 spot The wealth oil-to-telecom tycoon Mukesh Ambani also richest Indian surged rally shares Reliance Industries strong customer response Reliance Jio His younger brother Anil Ambani heads Reliance Communications.Here 10 things must know Forbes list Asia ''s richest families:1 The combined wealth Ambanis soared whopping 19 billion pipped South Korea ''s Lee family claim top notch rank Forbes list Asia ''s wealthiest clans.2 `` Because India ''s size Ambanis never dominant Samsung ''s Lees Korea But Mukesh ''s next generation establishing presence Reliance Jio telecom operation story could play years Tim Ferguson Editor Forbes Asia said Forbes website forbesindia.com.3 The Ambanis Indian family top 10 richest families list Asia `` No family highlights surge better Ambanis India year ''s biggest gainer dollar percentage terms Forbes said googletag.cmd.push function googletag.display `` adslotNativeVideo 4 Other richest Indian families list feature Premjis rank 11 19.2 billion net worth Hindujas 12th 18.8 billion Mittals 14th 17.2 billion Mistrys 16th 16.1 billion Birlas 19th 14.1 billion .5 Other Indian families got exalted list include Godrej family 20th 14 billion Bajajs 26th 9.3 billion Jindal family 32nd 7.68 billion Dabur ''s Burmans 35th 7.05 billion Eicher Motors Lals 36th 7 billion Shree Cement ''s Bangur family 37th 6.7 billion .6 Motherson Sumi Systems Sehgal family 41st 6.15 billion Wadia family 42nd 6.14 billion DLF ''s Kushal Pal Singh 44th 6.1 billion Patel family controls Cadila 45th 6 billion Piramals 47th 5.38 billion Munjals 48th 5.37 billion also made list.7 South Korea ''s Lee family ranked second Forbes list saw wealth soar 11.2 billion 40.8 billion year shares Samsung Electronics moved 75 per cent past year.8 Hong Kong ''s Kwok family controls Sun Hung Kai Properties ranked third year net wealth 40.4 billion Comments9 The list Asia ''s 50 Richest Families snapshot wealth using stock prices currency exchange rates close markets November 3 Private companies valued using financial ratios comparisons similar publicly traded firms.10 The ticket entry year ''s list 5 billion 1.6 billion 2016'
$ws.Range('C24').Value = 'news'
$ws.Range('B25').Value = 'Sania Mirza dropped top 10 first time since March 2014 © AFP Highlight Story widget starts Highlight Story widget end A consistent Divij Sharan broke top 50 career-high ranking Sania Mirza dropped top 10 first time since March 2014 per latest ATP WTA charts released Monday The left-handed Sharan reached final last three tournaments participated winning ATP 250 European Open ending runner-up Challenger events Brest Tashkent different partners He ranked exactly 50th gaining place The 31-year-old Delhi player let split long-time partner Purav Raja affect consistently well Tour despite stable partnership googletag.cmd.push function googletag.display ''adslotNativeVideo Now India two top 50 players doubles Rohan Bopanna 15 leading chart.Purav Raja 62 Leander Paes 70 Jeevan Nedunchezhiyan 97 top-100 Indian players `` made first ATP point 2003 It ''s long exciting journey since There lot hard-work sacrifices experiences travel decisions ups downs importantly grateful opportunity love playing tennis Sharan told PTI Bratislava Slovakia Sharan 19 ITF Futures titles 13 Challenger trophies three ATP Tour titles said `` belief win Grand Slam `` The guys losing close matches ones winning Slams We beaten sometimes feel win Grand Slam said `` employed Indian Oil since 2004 support massive career Head Adidas also providing racquets kit last four years added Sharan said next goal make top 30 play ATP 1000s Master series events `` It ''s good motivation encouragement Sharan reached pre-quarterfinals French Open year said.Also making good progress doubles Vishnu Vardhan jumped 16 places touch career-best ranking 118 Vardhan partner N Sriram Balaji 139 Shenzhen Challenger Saturday earn 90 ranking points each.In singles Yuki Bhambri 140 continues lead chart followed Ramkumar Ramanathan 148 Prajnesh Gunneswaran 255 Sumit Nagal 331 Balaji 350 .In WTA charts Ankita Raina India ''s top-ranked player number 281 followed Karman Kaur Thandi 307 Pranjala Yadlapalli 483 In doubles Sania placed number 12 following loss three places.Sania figured consistently top-10 since March 17 2013 The last time featured outside top-10 March 3 2013.The Hyderabadi world number one doubles entire 2016 season Her slide started season began number one could qualify year-end WTA Finale top-eight pairs top-eight singles players compete Comments Topics Sania Mirza Divij Sharan Tennis Get latest FIFA World Cup 2018 news check World Cup 2018 schedule Football live score World Cup standings Like us Facebook follow us Twitter FIFA World Cup 2018 updates You also download NDTV Cricket app Android iOS'
$ws.Range('C25').Value = 'sports'
$ws.Range('B26').Value = 'Mumbai The government working World Bank recognise 200 reforms help propel India top-50 bracket ease business according Department Industrial Policy Promotion DIPP secretary Ramesh Abhishek `` We already implemented 122 reforms year working World Bank recognise We initiate 90 ease business reforms year Mr Abhishek told reporters sidelines CII Invest North conclave.Yesterday India jumped 30 places rank 100th World Bank ''s ''ease business ranking helped slew reforms taxation licensing investor protection bankruptcy resolution `` A jump 30 positions World Bank rankings extremely commendable Now aim move top-50 countries bracket added.Mr Abhishek said department already started meeting stakeholders taking feedback reform measures taken government improve business climate `` This exercise helped us lot This time focused mainly taking feedback stakeholders All nodal ministries taking feedback said.googletag.cmd.push function googletag.display `` adslotNativeVideo The Secretary also said World Bank report acknowledged GST Goods Services Tax substantial reform `` hopefully positive impact country ''s ranking next year CommentsHe urged states address issues like inadequate infrastructure inflexible labour policies make competitive `` States invest skilling need come policies fixed term employment like government worked apparel sector We also working similar package leather footwear sector added Mr Abhishek'
$ws.Range('B27').Value = 'New Delhi Buoyed jump India ''s ease business ranking Finance Minister Arun Jaitley today vowed continue reforms help break top 50 nations.According World Bank India ''s rank ''ease business scale risen 130th 100th year helped slew reforms taxation licensing investor protection bankruptcy resolution.Addressing press conference minutes World Bank ranking made public Mr Jaitley said India major country named pursuing structural reforms `` In 2014 142nd last two years improved 131st 130th These generalised rankings It happened specific areas take tough parameters ranking said.The `` highest jump ranking possible significant improvements 10 judging parameters made last 3-4 years `` becomes easy business India said.The World Bank said named 10 countries structural reforms taken `` India one `` googletag.cmd.push function googletag.display `` adslotNativeVideo Mr Jaitley said 3 years ranking improved 142nd 100th sufficient work ongoing areas India lags behind `` There reason believe capacity significantly improve position said.The World Bank says India still lags areas starting business enforcing contracts dealing construction permits.It takes 30 days register new business 127 days 15 years ago `` number procedures still cumbersome local entrepreneurs still need go 12 procedures said `` The Prime Minister said target take India immediately first 50 believe doable therefore remaining 3-4 areas work done pushing greater force said.Elaborating ranking said significant improvement parameters protecting minority investors interest changes made company law `` We ranked 4th world This biggest jump field `` In availability credit India ranked 29th `` Biggest jump improvements done taxation front Till last year 189 nations 172nd It means paying taxes cumbersome process But year moved 53 places said.Also enactment new bankruptcy law India ''s ranking insolvency resolution also improved.The ranking include impact implementation Goods Services Tax GST created cut judging business environment June 1 `` In coming years expecting get credit said.Besides areas government done lot work work still progress get reflected numbers coming 1-2 years said.Giving examples said India ranked low 164th enforcement contracts arbitration act help improve this.Similarly states impressed upon give building construction permits online fast permissions handed `` One reform factored starting new business In 156th position created provision company act ... reform fully factored ranking sure future said.India ranked 156th latest ranking `` Work progress think able expedite work next months areas huge scope us jump 100 said.The finance minister said government seriously taking 3-4 issues India still lags behind next 1-2 years improvement 20-30 points ranking `` If see first available data August September think every data coming seems positive zone said CommentsMr Jaitley said reforms continuous process `` The 3-4 areas n''t treated high areas significant amount work progress impact shown ground ''m sure taken consideration next year ''s ranking said'
$ws.Range('B28').Value = 'Here ten points story Sharing details report Finance Minister Arun Jaitley said Tuesday India made biggest jump index `` The Prime Minister said top 50 The fact three years jumped 142 100 believe doable Mr Jaitley said.PM Modi series tweets said jump India ''s ranking result all-round multi-sectoral reform push `` Easier business environment leading historic opportunities entrepreneurs particularly MSME medium small scale sector sector bringing prosperity tweeted.The report based data New Delhi Mumbai ranked India among top 10 `` improvers globally done better eight 10 business indicators Mr Jaitley said government trying improve 10.India jumped 53 places 172 119 ease paying taxes leapt 4th rank protecting minority investors moved 15 places 44 29 ease businesses getting credit.There also red flags India fallen 16 places ease registering property slid slightly trading across borders one rank 155 156 ease starting businesses Mr Jaitley said would improve next year since report n''t factored changes made law Companies Act.New Zealand tops World Bank index ease business year Singapore number 2 Denmark placing third Somalia last China 78 Pakistan 147 Bangladesh 177 The report boost government defend economic policies amid criticism GDP growth rate hitting three-year low April June quarter.The government insisted fundamentals economy strong rejected allegations two big reforms last year ''s demonetisation new national tax GST launched year blame economic slowdown saying deliver big benefits `` Today ''s result clear signal India rest world country ready open business many decades competing preferred place business globally Annette Dixon World Bank ''s vice president South Asia told reporters New Delhi.The opposition Congress however kept attack saying NDA government ''s decisions like `` cease businesses medium small sectors The hype built government around report Congress said attempt seek `` manufactured relevance World Bank people suffering due GST notes ban.googletag.cmd.push function googletag.display `` adslotNativeVideo Comments'
$ws.Range('B29').Value = 'Mumbai October 17 Indian Institute Technology Bombay IIT Bombay moved one step 35th rank last year 34th rank year Quacquarelli Symonds QS Asia University Rankings Among top 50 Overall Institution improved 2016/17 performance one place Indian Institute Technology IIT Delhi ranked 41 IIT Madras placed rank 48 QS Asia University Rankings 2018 IISc Bangalore IIT Kanpur IIT Kharagpur University Delhi IIT Roorkee IIT Guwahati also found place 100 year.Earlier June year IIT Bombay ranked 179th QS World University Rankings showing improvement 2017 performance 40 places.The results Asia University Rankings release Monday October 16 2017 QS British company IIT Bombay ranks one India overall score 75.7 100 The Institute score 88.1 academic reputation 97.5 employer reputation 44.5 faculty student ration 96.3 staff PhD 89.6 papers per faculty 9.7 international faculty 6.3 international students 10.9 inbound exchange students 9.1 outbound exchange students 41.4 citations per paper Among ten parameters Employer Reputation indicated strongest one IIT Bombay rank 21 globally `` glad IIT Bombay continues rise rankings This reflection excellent work done faculty students Institute Director Instituet Prof. Devang Khakhar said.googletag.cmd.push function googletag.display `` adslotNativeVideo Considering approximately 11,900 universities globally made IIT Bombay one top 1.0 universities Asia.In QS Ranking 2018 released June year IIT Delhi dethroned IISc Bengaluru In 2017 rankings IISC placed highest among Indian institute managed rank top 200 IISc featured 152nd rank 2017 rankings QS Asia University Rankings 2018 Top 10Nanyang Technological University Singapore NTU Singapore Asia ''s top university according QS Asia University NTU renowned among world ''s best universities engineering technology rose two spots last year placed leading position QS Asia University league table ahead 400 universities Asia QS Asia University Rankings 2018 NTU Singapore ranked top list NTU made fantastic progress risen No 14 QS Asian rankings began 2009 The latest rankings outstanding achievement strong endorsement NTU ''s excellent reputation built years wish thank entire NTU community contributions NTU President Professor Bertil Anderson said.Universities Singapore Hong Kong China dominated rankings:1 Nanyang Technological University NTU Singapore 2 National University Singapore NUS 3 The Hong Kong University Science Technology Hong Kong 4 KAIST Korea Advanced Institute Science Technology South Korea 5 The University Hong Kong Hong Kong 6 Tsinghua University China 7 Fudan University China 8 City University Hong Kong Hong Kong 9 Peking University China 10 The Chinese University Hong Kong CUHK Click Education News Comments'
$ws.Range('B30').Value = 'Mumbai Banks may write 60 percent value bad loans 50 large stressed asset accounts take hit Rs 2.4 lakh crore ratings agency CRISIL said Wednesday.The 50 stressed companies account Rs 4 lakh crore soured loans largely metals construction power industries account nearly half total non-performing assets banking sector end March CRISIL said.The credit agency part S P Global estimated banks provisioned around 40 percent exposure assets CommentsBanks total non-performing loans Rs 7.29 lakh crore 5 percent India ''s gross domestic product end-March.The report comes rival India Ratings Research estimated banks taking 12 country ''s largest defaulters bankruptcy court would need make additional provisioning least Rs 18,000 crore.© Thomson Reuters 2017googletag.cmd.push function googletag.display `` adslotNativeVideo'
$ws.Range('B31').Value = 'New Delhi Enthused success SBI merger Finance Ministry considering clearing another proposal public sector banking space fiscal end goal create 4-5 global sized lenders.Five associate banks Bharatiya Mahila Bank became part SBI April 1 2017 catapulting country ''s largest lender among top 50 banks world.Now Finance Ministry looking replicate model case state-run banks reach critical mass compete global peers `` Consolidation must ... decision regard would based commercially prudent parameters If NPA situation gets better could one merger towards end fiscal senior official told PTI.Toxic loans public sector banks rose Rs 1 lakh crore Rs 6.06 lakh crore April-December 2016-17 bulk came power steel road infrastructure textile sectors.Finance Minister Arun Jaitley several occasions said India needs 5-6 banks global size scale consolidation banking sector done appropriate time.googletag.cmd.push function googletag.display `` adslotNativeVideo Whenever consolidation happens takes consideration interest stakeholders including employees shareholders official said adding win-win parties.A balancing act done merger given clearance various authorities regulators official added.The future merger proposals banking sector also require clearance Competition Commission India CCI see merged entity going create monopoly.In last consolidation drive saw light day April CCI nod needed case merger Bharatiya Mahila Bank BMB SBI There requirement merger associate banks SBI part parent.The Finance Ministry sought help government think-tank NITI Aayog global consultancy firms examine possibility next round consolidation PSU banks aim create lenders global size scale.NITI Aayog ''s report expected set tone roadmap consolidation future.There factors like regional balance geographical reach financial burden smooth human resource transition looked taking merger decision official said.The official added merger weak bank strong bank `` could pull latter `` There low-hanging fruits Big lenders like Bank Baroda take turnaround banks southern region Indian Overseas Bank Dena Bank could merged large South Indian bank official explained.Five associates BMB became part SBI April 1 2017 catapulting country ''s largest lender among top 50 banks world.State Bank Bikaner Jaipur SBBJ State Bank Hyderabad SBH State Bank Mysore SBM State Bank Patiala SBP State Bank Travancore SBT besides BMB merged SBI.With merger total customer base SBI reached around 37 crore branch network around 24,000 nearly 59,000 ATMs across country The merged entity began operation deposit base Rs 26 lakh crore advances Rs 18.50 lakh crore CommentsThe government February approved merger five associate banks SBI Later March cabinet approved merger BMB well SBI first merged State Bank Saurashtra 2008 Two years later State Bank Indore merged This story edited NDTV staff auto-generated syndicated feed'
